# dynamic_field.xlsx -- add 4 new "preferredLang" dynamic-field rows (100-103)
# and refresh the auto-filter / filter-database range to cover them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 100 ----------------------------------------------------------
$ws.Range("A100").Value = 10126
$ws.Range("B100").Value = "preferredLang"
$ws.Range("C100").Value = "user preferred Language"
$ws.Range("E100").Value = '{"value":"عربي","code":"ara"}'
$ws.Range("F100").Value = "eng"
$ws.Range("G100").Formula = "=TRUE()"
$ws.Range("H100").Value = "superadmin"
$ws.Range("I100").Value = "now()"

# ---- Row 101 ----------------------------------------------------------
$ws.Range("A101").Value = 10127
$ws.Range("B101").Value = "preferredLang"
$ws.Range("C101").Value = "Langue préférée de l'utilisateur"
$ws.Range("E101").Value = '{"value":"عربي","code":"ara"}'
$ws.Range("F101").Value = "fra"
$ws.Range("G101").Formula = "=TRUE()"
$ws.Range("H101").Value = "superadmin"
$ws.Range("I101").Value = "now()"

# ---- Row 102 ----------------------------------------------------------
$ws.Range("A102").Value = 10128
$ws.Range("B102").Value = "preferredLang"
$ws.Range("C102").Value = "user preferred Language"
$ws.Range("E102").Value = '{"value":"English","code":"eng"}'
$ws.Range("F102").Value = "ara"
$ws.Range("G102").Formula = "=TRUE()"
$ws.Range("H102").Value = "superadmin"
$ws.Range("I102").Value = "now()"

# ---- Row 103 ----------------------------------------------------------
$ws.Range("A103").Value = 10129
$ws.Range("B103").Value = "preferredLang"
$ws.Range("C103").Value = "user preferred Language"
$ws.Range("E103").Value = '{"value":"français","code":"fra"}'
$ws.Range("F103").Value = "ara"
$ws.Range("G103").Formula = "=TRUE()"
$ws.Range("H103").Value = "superadmin"
$ws.Range("I103").Value = "now()"

# ---- Formatting: wrap text on the text columns, boolean display format on G
foreach ($r in 100..103) {
    $ws.Range("B$r").WrapText = $true
    $ws.Range("C$r").WrapText = $true
    $ws.Range("E$r").WrapText = $true
    $ws.Range("F$r").WrapText = $true
    $ws.Range("H$r").WrapText = $true
    $ws.Range("I$r").WrapText = $true
    $ws.Range("G$r").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
    $ws.Rows.Item($r).RowHeight = 30
}

# ---- Refresh AutoFilter range to A1:I103 -------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:I103").AutoFilter()

# ---- Refresh the hidden _FilterDatabase defined name -------------------
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$I`$103"

# ---- Selection / active cell matches the saved view --------------------
$ws.Range("B99:B103").Select()
